$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new "TacX Flow refund" expense row (row 12) ---
# Copy formatting from the row above (row 11) so the new cells inherit the
# same number formats / styles used throughout the table.
$ws.Range("A11:H11").Copy()
$ws.Range("A12:H12").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A12").Value = 44481
$ws.Range("B12").Value = "TacX Flow refund"
$ws.Range("C12").Value = "Equipment"
$ws.Range("D12").Value = -200
$ws.Range("E12").Value = 1
$ws.Range("F12").Formula = "=E12*D12"
$ws.Range("H12").Value = "UPS refunded the seller because the package was damaged and missing components, she then gave me the price less shipping back."

# --- Fix up the running "cost per day" formulas (column G) ---
# They used to sum the per-item cost column (D); they should sum the
# extended sub-total column (F) instead, for every data row.
for ($i = 4; $i -le 47; $i++) {
    $cell = $ws.Range("G$i")
    $cell.Formula = $cell.Formula.Replace('$D$3:D', '$F$3:F')
}

# --- Autofit the columns that now have new content ---
$ws.Columns("D:D").AutoFit()
$ws.Columns("F:F").AutoFit()

# --- Update the view so the new row is visible/selected ---
[void]$ws.Range("E13").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3

$wb.Save()
